$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

$ws.Range("A1").Value = 'Name: uganda-water-source-comm-func'
$ws.Range("A2").Value = 'Description: Functional Water Source Committees (%)'
$ws.Range("A3").Value = 'Units of measure: %'
$ws.Range("A4").Value = 'Source: Water and Environment Sector Performance Reports 2010-2014 - Ministry of Water and Environment'
$ws.Range("A5").Value = 'Source-link: http://www.mwe.go.ug/index.php?option=com_docman&task=cat_view&Itemid=223&gid=15'
$ws.Range("A6").Value = ''
$ws.Range("A7").Value = 'Notes:'
$ws.Range("A8").Value = ''
$ws.Range("A9").Value = 'On the ''Data-wide-value'' sheet, we have provided the indicator in a wide format. The values you see listed there are from the ''value'' column.'
$ws.Range("A10").Value = ''
$ws.Range("A11").Value = ''
$ws.Range("A12").Value = ''
$ws.Range("A13").Value = 'The following is data downloaded from Development Initiative''s Datahub: http://devinit.org/data'
$ws.Range("A14").Value = 'It is licensed under a Creative Commons Attribution 4.0 International license.'
$ws.Range("A15").Value = 'More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/'
$ws.Range("A16").Value = 'For concerns, questions, or corrections: please email info@devinit.org'
$ws.Range("A17").Value = 'Copyright Development Initiatives Poverty Research Ltd. 2015'
